$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 952
$ws.Range("I2").Value = 142.85715
$ws.Range("K2").Value = 142.85715
$ws.Range("M2").Value = -29.85714999999999
$ws.Range("H9").Value = 415.33334
$ws.Range("I9").Value = 431.1
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 431.1
$ws.Range("L9").Value = 100
$ws.Range("M9").Value = -262.1
$ws.Range("N9").Value = -438
$ws.Range("H17").Value = 1561.439
$ws.Range("J17").Value = 1561.439
$ws.Range("L17").Value = 4684.317
$ws.Range("N17").Value = -5020.317
$ws.Range("H43").Value = 2750
$ws.Range("I43").Value = 2333.3333
$ws.Range("K43").Value = 2333.3333
$ws.Range("M43").Value = -2264.3333
$ws.Range("H51").Value = 6077.143
$ws.Range("I51").Value = 7660
$ws.Range("J51").Value = 3966.6667
$ws.Range("K51").Value = 7660
$ws.Range("L51").Value = 3966.6667
$ws.Range("M51").Value = -7176
$ws.Range("N51").Value = -4934.6667
$ws.Range("H100").Value = 8027.5264
$ws.Range("I100").Value = 1028.5
$ws.Range("J100").Value = 13117.728
$ws.Range("K100").Value = 1028.5
$ws.Range("L100").Value = 13117.728
$ws.Range("M100").Value = -487.5
$ws.Range("N100").Value = -14199.728
$ws.Range("H113").Value = 5549.3
$ws.Range("I113").Value = 3998.8
$ws.Range("J113").Value = 7099.8
$ws.Range("K113").Value = 3998.8
$ws.Range("L113").Value = 7099.8
$ws.Range("M113").Value = -744.8000000000002
$ws.Range("N113").Value = -13607.8
$ws.Range("H116").Value = 7249.8335
$ws.Range("I116").Value = 6700
$ws.Range("J116").Value = 9999
$ws.Range("K116").Value = 6700
$ws.Range("L116").Value = 9999
$ws.Range("M116").Value = -3258
$ws.Range("N116").Value = -16883
$ws.Range("H132").Value = 2252.2112
$ws.Range("I132").Value = 1817.2881
$ws.Range("K132").Value = 5451.8643
$ws.Range("M132").Value = -2921.8643
$ws.Range("H134").Value = 63332.8
$ws.Range("J134").Value = 63332.8
$ws.Range("L134").Value = 63332.8
$ws.Range("N134").Value = -73472.8
$ws.Range("H135").Value = 556520.75
$ws.Range("I135").Value = 606855.9399999999
$ws.Range("J135").Value = 2833.3333
$ws.Range("K135").Value = 5461703.459999999
$ws.Range("L135").Value = 25499.9997
$ws.Range("M135").Value = -5459168.459999999
$ws.Range("N135").Value = -30569.9997
$ws.Range("H138").Value = 4401.109
$ws.Range("I138").Value = 2631.9524
$ws.Range("J138").Value = 5493.8237
$ws.Range("K138").Value = 7895.8572
$ws.Range("L138").Value = 16481.4711
$ws.Range("M138").Value = -2755.8572
$ws.Range("N138").Value = -26761.4711
$ws.Range("H141").Value = 1893.1029
$ws.Range("I141").Value = 844.87933
$ws.Range("K141").Value = 2534.63799
$ws.Range("M141").Value = 2645.36201

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H45").Value = 2969
$ws.Range("I45").Value = 2486.7
$ws.Range("K45").Value = 2486.7
$ws.Range("M45").Value = -2109.7
$ws.Range("H74").Value = 1966.3334
$ws.Range("I74").Value = 1876.1177
$ws.Range("J74").Value = 3500
$ws.Range("K74").Value = 1876.1177
$ws.Range("L74").Value = 3500
$ws.Range("M74").Value = -1002.1177
$ws.Range("N74").Value = -5248
$ws.Range("H77").Value = 1966.3334
$ws.Range("I77").Value = 1876.1177
$ws.Range("J77").Value = 3500
$ws.Range("K77").Value = 9380.5885
$ws.Range("L77").Value = 17500
$ws.Range("M77").Value = -5012.5885
$ws.Range("N77").Value = -26236
$ws.Range("H102").Value = 2975.8
$ws.Range("I102").Value = 3399.5
$ws.Range("K102").Value = 3399.5
$ws.Range("M102").Value = -1777.5
$ws.Range("H122").Value = 3310.95
$ws.Range("I122").Value = 1801.1154
$ws.Range("K122").Value = 5403.3462
$ws.Range("M122").Value = -2953.3462
$ws.Range("H132").Value = 5132.483
$ws.Range("I132").Value = 3366.3333
$ws.Range("K132").Value = 10098.9999
$ws.Range("M132").Value = -7568.999899999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 14079.299
$ws.Range("I134").Value = 1675.303
$ws.Range("K134").Value = 5025.909000000001
$ws.Range("M134").Value = -2490.909000000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 308175.38
$ws.Range("I31").Value = 334959.97
$ws.Range("J31").Value = 174252.33
$ws.Range("K31").Value = 334959.97
$ws.Range("L31").Value = 174252.33
$ws.Range("M31").Value = -334664.97
$ws.Range("N31").Value = -174842.33
$ws.Range("H34").Value = 308175.38
$ws.Range("I34").Value = 334959.97
$ws.Range("J34").Value = 174252.33
$ws.Range("K34").Value = 334959.97
$ws.Range("L34").Value = 174252.33
$ws.Range("M34").Value = -334757.97
$ws.Range("N34").Value = -174656.33
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H51").Value = 55074.25
$ws.Range("J51").Value = 55074.25
$ws.Range("L51").Value = 55074.25
$ws.Range("N51").Value = -56546.25
$ws.Range("H52").Value = 70077.71000000001
$ws.Range("J52").Value = 69818.8
$ws.Range("L52").Value = 69818.8
$ws.Range("N52").Value = -70406.8
$ws.Range("H58").Value = 217126.52
$ws.Range("I58").Value = 347541.62
$ws.Range("J58").Value = 7013.278
$ws.Range("K58").Value = 347541.62
$ws.Range("L58").Value = 7013.278
$ws.Range("M58").Value = -347338.62
$ws.Range("N58").Value = -7419.278
$ws.Range("H60").Value = 19666
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H61").Value = 55074.25
$ws.Range("J61").Value = 55074.25
$ws.Range("L61").Value = 55074.25
$ws.Range("N61").Value = -55770.25
$ws.Range("H99").Value = 5998.7085
$ws.Range("I99").Value = 5215.643
$ws.Range("K99").Value = 5215.643
$ws.Range("M99").Value = -3717.643
$ws.Range("H126").Value = 5998.7085
$ws.Range("I126").Value = 5215.643
$ws.Range("K126").Value = 15646.929
$ws.Range("M126").Value = -13176.929
$ws.Range("H136").Value = 217126.52
$ws.Range("I136").Value = 347541.62
$ws.Range("J136").Value = 7013.278
$ws.Range("K136").Value = 1042624.86
$ws.Range("L136").Value = 21039.834
$ws.Range("M136").Value = -1040074.86
$ws.Range("N136").Value = -26139.834

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 47.07143
$ws.Range("I6").Value = 47.07143
$ws.Range("K6").Value = 141.21429
$ws.Range("M6").Value = -28.21429000000001
$ws.Range("H136").Value = 8161.222
$ws.Range("I136").Value = 6245.5
$ws.Range("J136").Value = 11992.667
$ws.Range("K136").Value = 18736.5
$ws.Range("L136").Value = 35978.001
$ws.Range("M136").Value = -13636.5
$ws.Range("N136").Value = -46178.001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 142867140
$ws.Range("J29").Value = 166674990
$ws.Range("L29").Value = 166674990
$ws.Range("N29").Value = -166675570

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 917969.6
$ws.Range("I7").Value = 30004
$ws.Range("K7").Value = 30004
$ws.Range("M7").Value = -29892
$ws.Range("H16").Value = 1456.625
$ws.Range("I16").Value = 1333.0769
$ws.Range("J16").Value = 1992
$ws.Range("K16").Value = 1333.0769
$ws.Range("L16").Value = 1992
$ws.Range("M16").Value = -1163.0769
$ws.Range("N16").Value = -2332
$ws.Range("H31").Value = 1919.1666
$ws.Range("I31").Value = 1303
$ws.Range("J31").Value = 5000
$ws.Range("K31").Value = 1303
$ws.Range("L31").Value = 5000
$ws.Range("M31").Value = -1055
$ws.Range("N31").Value = -5496
$ws.Range("H93").Value = 3380.889
$ws.Range("I93").Value = 3055.8
$ws.Range("J93").Value = 3787.25
$ws.Range("K93").Value = 3055.8
$ws.Range("L93").Value = 3787.25
$ws.Range("M93").Value = -1807.8
$ws.Range("N93").Value = -6283.25
$ws.Range("H126").Value = 917969.6
$ws.Range("I126").Value = 30004
$ws.Range("K126").Value = 90012
$ws.Range("M126").Value = -87542
$ws.Range("H132").Value = 5129.771
$ws.Range("I132").Value = 4639.0884
$ws.Range("K132").Value = 13917.2652
$ws.Range("M132").Value = -11387.2652
$ws.Range("H134").Value = 42899.7
$ws.Range("J134").Value = 42899.7
$ws.Range("L134").Value = 42899.7
$ws.Range("N134").Value = -53039.7

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H107").Value = 977.6667
$ws.Range("I107").Value = 1389.5
$ws.Range("K107").Value = 4168.5
$ws.Range("M107").Value = -2248.5
$ws.Range("H122").Value = 45459584
$ws.Range("I122").Value = 142859780
$ws.Range("K122").Value = 428579340
$ws.Range("M122").Value = -428576890
$ws.Range("H132").Value = 21567.623
$ws.Range("J132").Value = 76285.36
$ws.Range("L132").Value = 228856.08
$ws.Range("N132").Value = -233916.08
$ws.Range("H136").Value = 543098
$ws.Range("I136").Value = 1158339.5
$ws.Range("J136").Value = 147585.58
$ws.Range("K136").Value = 3475018.5
$ws.Range("L136").Value = 442756.74
$ws.Range("M136").Value = -3472468.5
$ws.Range("N136").Value = -447856.74
